$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in the title row
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 08:22"

# Irlanda (row 23): refresh Casos activos / Recuperados
$ws.Range("D23").Value = 77
$ws.Range("E23").Value = 10996

# Moldavia (row 59): refresh Casos activos / Recuperados / Muertes
$ws.Range("D59").Value = 171
$ws.Range("E59").Value = 1722
$ws.Range("H59").Value = 41

# Taiwan & Malta swap places (row 102/103) with Taiwan's data refreshed
$ws.Range("A102").Value = "Taiwan"
$ws.Range("B102").Value = 395
$ws.Range("C102").Value = 2
$ws.Range("D102").Value = 137
$ws.Range("E102").Value = 252
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 6

$ws.Range("A103").Value = "Malta"
$ws.Range("B103").Value = 393
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 44
$ws.Range("E103").Value = 346
$ws.Range("F103").Value = 4
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 3

# El Salvador & Martinica swap places (row 124/125) with El Salvador's data refreshed
$ws.Range("A124").Value = "El Salvador"
$ws.Range("B124").Value = 159
$ws.Range("C124").Value = 10
$ws.Range("D124").Value = 30
$ws.Range("E124").Value = 123
$ws.Range("F124").Value = 2
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 6

$ws.Range("A125").Value = "Martinica"
$ws.Range("B125").Value = 158
$ws.Range("C125").Value = 1
$ws.Range("D125").Value = 73
$ws.Range("E125").Value = 77
$ws.Range("F125").Value = 17
$ws.Range("G125").Value = 2
$ws.Range("H125").Value = 8
